$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 1 was "Adding weights to each objective"
# Paragraph 2 was "Increasing number of objectives"
# They get swapped, and each gets split into multiple runs.

$para1 = $tr.Paragraphs(1)
$para1.Text = "Increasing "
$para1.InsertAfter("number of ") | Out-Null
$para1.InsertAfter("objectives") | Out-Null

$para2 = $tr.Paragraphs(2)
$para2.Text = "Adding weights to "
$para2.InsertAfter("each ") | Out-Null
$para2.InsertAfter("objective") | Out-Null
